# Add more blank ("left blank for whiteboard") slides to the default
# presentation. The deck already ends with a "This slide left blank for
# whiteboard" slide (the last slide); we clone that slide 14 more times
# so the deck grows from 6 slides to 20, all the new ones appended at
# the end in order, each an exact duplicate (same layout/shapes/text)
# of the existing blank-whiteboard slide.

$p = $ppt.ActivePresentation

$blankSlideCount = 14

for ($i = 0; $i -lt $blankSlideCount; $i++) {
    $lastIndex = $p.Slides.Count
    $source = $p.Slides.Item($lastIndex)
    $source.Duplicate() | Out-Null
}
